# Applies the "blockTemplate" sheet addition + openLeftMenu() parameterization
# described by the commit "Generating and pushing into server. Also deploys functions".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) tools sheet: addBlock / editBlock now call openLeftMenu() with an argument
# ---------------------------------------------------------------------------
$tools = $wb.Worksheets.Item("tools")
$tools.Range("C5").Value2 = "openLeftMenu('ADD_BLOCK')"
$tools.Range("C6").Value2 = "openLeftMenu('EDIT_BLOCK')"

# ---------------------------------------------------------------------------
# 2) Add the new "blockTemplate" worksheet after the last sheet (networkIcon)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bt = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$bt.Name = "blockTemplate"

# Header row
$bt.Range("A1").Value2 = "name"
$bt.Range("B1").Value2 = "type"
$bt.Range("C1").Value2 = "argument"
$bt.Range("E1").Formula = '="const "&MID(CELL("filename",A1),FIND("]",CELL("filename",A1))+1,255)&"=["&_xlfn.CONCAT(D:D)&"];"'

# Field definitions (name / type / argument)
$rows = @(
    @("Category",        "string",   ""),
    @("Name",             "string",   ""),
    @("Description",      "string",   ""),
    @("Parameters",       "json",     ""),
    @("Label",            "function", ""),
    @("MaxInTerminals",   "float",    ""),
    @("MaxOutTerminals",  "float",    ""),
    @("Icon",             "function", ""),
    @("Init",             "function", ""),
    @("End",              "function", ""),
    @("Constructor",      "function", "Data"),
    @("Destructor",       "function", "Data"),
    @("RunTimeExec",      "function", ""),
    @("Evaluate",         "function", ""),
    @("Details",          "function", ""),
    @("ValidateParams",   "function", "")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $bt.Cells.Item($r, 1).Value2 = $rows[$i][0]
    $bt.Cells.Item($r, 2).Value2 = $rows[$i][1]
    if ($rows[$i][2] -ne "") {
        $bt.Cells.Item($r, 3).Value2 = $rows[$i][2]
    }
}

# D2 is entered on its own, then D3:D17 is filled as one shared-formula block
$bt.Range("D2").Formula = '="{"&CHAR(34)&$A$1&CHAR(34)&":"&CHAR(34)&A2&CHAR(34)&","&CHAR(34)&$B$1&CHAR(34)&":"&CHAR(34)&B2&CHAR(34)&","&CHAR(34)&$C$1&CHAR(34)&":"&CHAR(34)&C2&CHAR(34)&"}"&IF(ISBLANK(A3),"",",")'
$bt.Range("D3:D17").Formula = '="{"&CHAR(34)&$A$1&CHAR(34)&":"&CHAR(34)&A3&CHAR(34)&","&CHAR(34)&$B$1&CHAR(34)&":"&CHAR(34)&B3&CHAR(34)&","&CHAR(34)&$C$1&CHAR(34)&":"&CHAR(34)&C3&CHAR(34)&"}"&IF(ISBLANK(A4),"",",")'

# I2 is entered on its own, then I3:I17 is filled as one shared-formula block
$bt.Range("I2").Formula = "=TRIM(A2)"
$bt.Range("I3:I17").Formula = "=TRIM(A3)"

# Column A width, matching the authored sheet (~16 chars, best-fit)
$bt.Columns.Item(1).AutoFit() | Out-Null
$bt.Columns.Item(1).ColumnWidth = 15.1666666

# Selection left on A5 for the new sheet (not the active sheet)
$bt.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) View bookkeeping: "tools" becomes the active tab, "networkIcon" keeps
#    C1 selected but is no longer the active tab.
# ---------------------------------------------------------------------------
$networkIcon = $wb.Worksheets.Item("networkIcon")
$networkIcon.Activate()
$networkIcon.Range("C1").Select() | Out-Null

$tools.Activate()
$tools.Range("J1").Select() | Out-Null
